$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster -> ECs, plus updated TPM-derived metrics)
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.072366
$ws.Range("H2").Value = 0.217098
$ws.Range("M2").Value = 62.87391666666667
$ws.Range("N2").Value = 188.62175
$ws.Range("O2").Value = 0.5484251561826182
$ws.Range("P2").Value = 0.5484251561826182
$ws.Range("Q2").Value = 4.5499338535
$ws.Range("R2").Value = 40.9494046815
$ws.Range("S2").Value = 0.5484251561826182
$ws.Range("T2").Value = 0.5484251561826182

# Row 3 (Target cluster -> FAPs, plus updated TPM-derived metrics)
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.072366
$ws.Range("H3").Value = 0.217098
$ws.Range("O3").Value = 0.3074331251635
$ws.Range("P3").Value = 0.3074331251635
$ws.Range("Q3").Value = 2.550576624904
$ws.Range("R3").Value = 22.955189624136
$ws.Range("S3").Value = 0.3074331251635
$ws.Range("T3").Value = 0.3074331251635

# Row 4 (Target cluster -> MuSCs, plus updated TPM-derived metrics)
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.072366
$ws.Range("H4").Value = 0.217098
$ws.Range("M4").Value = 7.926563000000001
$ws.Range("N4").Value = 23.779689
$ws.Range("O4").Value = 0.06914038096772557
$ws.Range("P4").Value = 0.06914038096772555
$ws.Range("Q4").Value = 0.573613658058
$ws.Range("R4").Value = 5.162522922522
$ws.Range("S4").Value = 0.06914038096772557
$ws.Range("T4").Value = 0.06914038096772555

# Row 5 (Target cluster stays Resolving-Mac, but updated TPM-derived metrics)
$ws.Range("G5").Value = 0.072366
$ws.Range("H5").Value = 0.217098
$ws.Range("M5").Value = 8.598489333333333
$ws.Range("N5").Value = 25.795468
$ws.Range("O5").Value = 0.07500133768615619
$ws.Range("P5").Value = 0.07500133768615617
$ws.Range("Q5").Value = 0.6222382790959999
$ws.Range("R5").Value = 5.600144511864
$ws.Range("S5").Value = 0.07500133768615619
$ws.Range("T5").Value = 0.07500133768615617
